# Insert a new column before column A to hold row identifiers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Columns.Item(1).Insert()

# Match the header formatting used by the other header cells (bold, thin border, centered)
$ws.Range("B1").Copy() | Out-Null
$ws.Range("A1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Header for the new ID column
$ws.Range("A1").Value = "ID"

# Row identifiers for rows 2-25
$ids = @("Hb 2","Hb 3","S 24","S 28","Hb 107","Hb 66","Hb 69","Hb 95","Hb 99","Hb 92","Hb 40","Hb 41","S 11","Hb 57","S 21","S 22","S 3","S 4","S 5","Hb 74","Hb 79","Hb 32","S 15","S 16")

for ($i = 0; $i -lt $ids.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $ids[$i]
}
